$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 52; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    $cell.Value = $current + 1
}
